# The workbook has one worksheet per year (tabs "2000" .. "2100"). Each
# sheet has an identical small table where row 6 is labelled "Onshore wind
# plants" and row 7 is labelled "Offshore wind plants" (column C), with the
# numeric series living in column E (D/F/G are always 0 for these rows).
#
# The commit swaps the Onshore/Offshore rows on every sheet: the label in
# row 6 becomes "Offshore wind plants" (and row 7 becomes "Onshore wind
# plants"), while the E-column figures move with their original label, i.e.
# E6 and E7 swap values.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Swap the numeric values between row 6 and row 7 (column E).
    $onshoreValue = $ws.Range("E6").Value2
    $offshoreValue = $ws.Range("E7").Value2
    $ws.Range("E6").Value2 = $offshoreValue
    $ws.Range("E7").Value2 = $onshoreValue

    # Swap the row labels so row 6 is now "Offshore" and row 7 "Onshore".
    $ws.Range("C6").Value2 = "Offshore wind plants"
    $ws.Range("C7").Value2 = "Onshore wind plants"
}
